# Insert a new "REGULATORY_EXPIRE_DATE" column before the existing
# "ROUTE_OF_ADMINISTRATION" column (Q), shifting Q..AA to R..AB.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1").EntireColumn.Insert()
$ws.Range("Q1").Value = "REGULATORY_EXPIRE_DATE"
